$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.563.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.787.82'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '351.41'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.46'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.549'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.91%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.623'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.46'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0833'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.95'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.77'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.222.13'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.802.21'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.928'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.563.10'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.70'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.12'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.35'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0967'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.57'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.54'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.76'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.164'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.28'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.59%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +9.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '51.89'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.66'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +8.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0442'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -6.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0853'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.56'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.12'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.74%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '120.11'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.99'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.19'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.134.04'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.35'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.93%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.46%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +17.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.45'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.09%  '
$ws.Range("B51").Value = 'SEI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.896'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -5.55%  '
